# ------------------------------------------------------------------
# Append new scrape run (2025-10-28 01:17:29 JST):
#   - a brand-new listing ("UberEats...") is inserted as row 10
#   - every previously-existing listing (old rows 10-20) shifts down by
#     one row (new rows 11-21)
#   - rows 2-9 keep their data but get the refreshed "取得日時" timestamp
#   - column F hyperlinks are rebuilt so they keep tracking their own URL
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$timestamp = '2025-10-28 01:17:29'

# Final state for rows 2..21, columns A..H (taken straight from the diff):
$rows = @(
    @('2025-10-28 01:17:29', '医療機関向けAIアプリとLINEの連携開発を支援してくださるAIエンジニア募集(AI/バックエンド)', 'システム開発', '300,000 円 ~ 500,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5416301', 385, '🔥AI,Ai ◆開発 ◇アプリ'),
    @('2025-10-28 01:17:29', '大企業の業務効率化AIプロジェクトの技術方針策定を支援するAIテックリード募集', 'システム開発', '100,000 円 ~ 200,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5416307', 378, '🔥AI,Ai ◆効率化'),
    @('2025-10-28 01:17:29', 'Azureでの社内文書検索RAG開発の精度改善を伴走支援してくださるAIエンジニア募集', 'システム開発', '300,000 円 ~ 500,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5416305', 375, '🔥AI,Ai ◆開発'),
    @('2025-10-28 01:17:29', 'Stable Diffusionに詳しいLoRAなどを用いた画像生成AIエンジニア募集', 'システム開発', '300,000 円 ~ 500,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5416328', 310, '🔥AI,Ai'),
    @('2025-10-28 01:17:29', 'Webシステム チャット機能へのChatwork連携API新規開発・組み込み', 'システム開発', '300,000 円 ~ 500,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5421556', 265, '🔥API ◆開発'),
    @('2025-10-28 01:17:29', '【急募】PDF見積書をExcel注文書に変換するシステム開発', 'システム開発', '100,000 円 ~ 200,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5421561', 118, '◆開発,システム開発'),
    @('2025-10-28 01:17:29', '【Unity/XRエンジニア募集】製造業DX支援!既存システムと連携するXRアプリ開発', 'システム開発', '300,000 円 ~ 500,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5421687', 115, '◆開発 ◇アプリ'),
    @('2025-10-28 01:17:29', 'Amazon購入履歴の明細PDFを自動ダウンロード&自動リネームするシステム開発', 'システム開発', '20,000 円 ~ 50,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5421083', 113, '◆開発,システム開発'),
    @('2025-10-28 01:17:29', '【急募】UberEats案件オファー抽出アプリのバックエンド開発', 'システム開発', '100,000 円 ~ 200,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5421779', 93, '◆開発 ◇アプリ'),
    @('2025-10-28 01:17:29', '【RPA構築依頼】不動産問い合わせ対応自動化(アシロボ使用、報酬10万円)', 'システム開発', '100,000 円 ~ 200,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5421443', 88, '◆自動化'),
    @('2025-10-28 01:17:29', '初回 楽天RMSの配布型クーポンの自動登録システムの開発', 'システム開発', '20,000 円 ~ 50,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5421265', 78, '◆開発'),
    @('2025-10-28 01:17:29', '【カフェ情報プラットフォーム開発】基本設計からリリースまで', 'システム開発', '300,000 円 ~ 500,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5420868', 75, '◆開発'),
    @('2025-10-28 01:17:29', '【簡単RPA構築】特定ツールからのデータ取得・Excel処理・スプレッドシートへの貼付け', 'システム開発', '20,000 円 ~ 50,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5421445', 68, '◆ツール'),
    @('2025-10-28 01:17:29', '音声デシベル検知器の開発を手伝ってくれる方募集!', 'システム開発', '200,000 円 ~ 300,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5421105', 68, '◆開発'),
    @('2025-10-28 01:17:29', 'WordPressサイトのリニューアル作業依頼', 'システム開発', '10,000 円 ~ 20,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5420971', 50, '◇サイト ○WordPress'),
    @('2025-10-28 01:17:29', '【カンタン作業】サイト環境立ち上げ検証の作業!', 'システム開発', '~ 5,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5421230', 30, '◇サイト'),
    @('2025-10-28 01:17:29', '〖リモート可〗Delphiエンジニア募集', 'システム開発', '300,000 円 ~ 500,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5341051', 25, $null),
    @('2025-10-28 01:17:29', '【急募】Shopifyでのフォーム一体型LPコード作成依頼', 'システム開発', '20,000 円 ~ 50,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5421564', 13, $null),
    @('2025-10-28 01:17:29', '【急募】年末調整業務のマクロ作成依頼', 'システム開発', '20,000 円 ~ 50,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5421418', 13, $null),
    @('2025-10-28 01:17:29', '【10,000円1万枚】指定したURL先のHPのキャプチャー画像を作成お願い致します。', 'システム開発', '5,000 円 ~ 10,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5421177', 10, $null)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $timestamp
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    if ($data[7] -eq $null) {
        $ws.Cells.Item($r, 8).Value = ""
    } else {
        $ws.Cells.Item($r, 8).Value = $data[7]
    }
}

# Rebuild the column-F hyperlinks so each one keeps pointing at the URL
# that is now actually displayed in that row (Hyperlinks.Delete() clears the
# whole collection in this engine, so everything is re-added in final order).
$ws.Hyperlinks.Delete()
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $url = $rows[$i][5]
    $ws.Hyperlinks.Add($ws.Range("F" + $r), $url)
}

Write-Output "done"
